$d = $word.ActiveDocument

$d.Content.Find.Execute("754×5=", $true, $false, $false, $false, $false, $true, 1, $false, "396×3=", 2)
$d.Content.Find.Execute("863×9=", $true, $false, $false, $false, $false, $true, 1, $false, "750×3=", 2)
$d.Content.Find.Execute("693×9=", $true, $false, $false, $false, $false, $true, 1, $false, "337×5=", 2)
$d.Content.Find.Execute("802×2=", $true, $false, $false, $false, $false, $true, 1, $false, "686×9=", 2)
$d.Content.Find.Execute("545×2=", $true, $false, $false, $false, $false, $true, 1, $false, "406×5=", 2)
$d.Content.Find.Execute("531×4=", $true, $false, $false, $false, $false, $true, 1, $false, "584×5=", 2)
$d.Content.Find.Execute("478×5=", $true, $false, $false, $false, $false, $true, 1, $false, "445×3=", 2)
$d.Content.Find.Execute("436×7=", $true, $false, $false, $false, $false, $true, 1, $false, "748×4=", 2)
$d.Content.Find.Execute("319×9=", $true, $false, $false, $false, $false, $true, 1, $false, "696×4=", 2)
$d.Content.Find.Execute("952×6=", $true, $false, $false, $false, $false, $true, 1, $false, "255×9=", 2)
$d.Content.Find.Execute("430×7=", $true, $false, $false, $false, $false, $true, 1, $false, "926×8=", 2)
$d.Content.Find.Execute("898×8=", $true, $false, $false, $false, $false, $true, 1, $false, "496×8=", 2)
$d.Content.Find.Execute("643×5=", $true, $false, $false, $false, $false, $true, 1, $false, "488×8=", 2)
$d.Content.Find.Execute("288×7=", $true, $false, $false, $false, $false, $true, 1, $false, "252×9=", 2)
$d.Content.Find.Execute("343×9=", $true, $false, $false, $false, $false, $true, 1, $false, "759×7=", 2)
$d.Content.Find.Execute("408×5=", $true, $false, $false, $false, $false, $true, 1, $false, "406×2=", 2)
$d.Content.Find.Execute("254×7=", $true, $false, $false, $false, $false, $true, 1, $false, "964×6=", 2)
$d.Content.Find.Execute("108×7=", $true, $false, $false, $false, $false, $true, 1, $false, "374×2=", 2)
$d.Content.Find.Execute("640×7=", $true, $false, $false, $false, $false, $true, 1, $false, "379×4=", 2)
$d.Content.Find.Execute("675×5=", $true, $false, $false, $false, $false, $true, 1, $false, "325×9=", 2)
$d.Content.Find.Execute("823×5=", $true, $false, $false, $false, $false, $true, 1, $false, "432×4=", 2)
$d.Content.Find.Execute("549×6=", $true, $false, $false, $false, $false, $true, 1, $false, "150×6=", 2)
$d.Content.Find.Execute("256×7=", $true, $false, $false, $false, $false, $true, 1, $false, "977×7=", 2)
$d.Content.Find.Execute("236×2=", $true, $false, $false, $false, $false, $true, 1, $false, "793×3=", 2)
$d.Content.Find.Execute("597×2=", $true, $false, $false, $false, $false, $true, 1, $false, "358×2=", 2)
